$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the "last updated" timestamp shown in A1
#    "Datos actualizados a 30 de Marzo de 2020 a las 08:50"
#      -> "Datos actualizados a 30 de Marzo de 2020 a las 09:20"
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "Datos actualizados a 30 de Marzo de 2020 a las 09:20"

# ---------------------------------------------------------------------------
# 2) Re-order a few countries (their position in the list changes, which
#    shifts which stats line up with which country label) and refresh the
#    case counts for the affected rows with the newer snapshot values.
# ---------------------------------------------------------------------------

# Country labels (column A) that move to a new row
$countryUpdates = @{
    20 = "Israel"
    21 = "Noruega"
    22 = "Brasil"
    32 = "Polonia"
    33 = "Japon"
    66 = "Lituania"
    67 = "Marruecos"
    68 = "Ucrania"
}

foreach ($row in $countryUpdates.Keys) {
    $ws.Range("A$row").Value2 = $countryUpdates[$row]
}

# Updated statistics (columns B..H = Casos totales, Nuevos casos, Casos
# activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$statUpdates = @{
    4   = @(142735, 275, 4562, 135684, 2970, 5, 2489)
    17  = @(8867, 79, 479, 8302, 187, 0, 86)
    20  = @(4347, 100, 132, 4200, 80, 0, 15)
    21  = @(4305, 21, 7, 4272, 91, 0, 26)
    22  = @(4256, 0, 6, 4114, 296, 0, 136)
    32  = @(1905, 43, 7, 1872, 3, 4, 26)
    33  = @(1866, 0, 424, 1388, 56, 0, 54)
    66  = @(484, 24, 1, 476, 2, 0, 7)
    67  = @(479, 0, 13, 440, 1, 0, 26)
    68  = @(475, 0, 6, 459, 0, 0, 10)
    82  = @(263, 0, 13, 248, 33, 0, 2)
    116 = @(91, 6, 0, 90, 1, 0, 1)
    126 = @(49, 1, 19, 25, 1, 0, 5)
}

$cols = @("B", "C", "D", "E", "F", "G", "H")

foreach ($row in $statUpdates.Keys) {
    $values = $statUpdates[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value2 = $values[$i]
    }
}
